# Auto-generated edit script applying updated leve market-price data
# across all 8 crafting-class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 167.5
$ws.Range("I5").Value = 185.71428
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 185.71428
$ws.Range("L5").Value = 40
$ws.Range("M5").Value = -70.71428
$ws.Range("N5").Value = -270
$ws.Range("H51").Value = 1738.7142
$ws.Range("J51").Value = 1720.8
$ws.Range("L51").Value = 1720.8
$ws.Range("N51").Value = -2688.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 521
$ws.Range("I3").Value = 521
$ws.Range("K3").Value = 521
$ws.Range("M3").Value = -406
$ws.Range("H8").Value = 43337.332
$ws.Range("J8").Value = 43337.332
$ws.Range("L8").Value = 43337.332
$ws.Range("N8").Value = -43625.332
$ws.Range("H44").Value = 30449
$ws.Range("J44").Value = 30449
$ws.Range("L44").Value = 30449
$ws.Range("N44").Value = -31425
$ws.Range("H55").Value = 23735.334
$ws.Range("J55").Value = 23735.334
$ws.Range("L55").Value = 23735.334
$ws.Range("N55").Value = -24365.334
$ws.Range("H80").Value = 33402.855
$ws.Range("J80").Value = 33402.855
$ws.Range("L80").Value = 33402.855
$ws.Range("N80").Value = -35398.855
$ws.Range("H83").Value = 33402.855
$ws.Range("J83").Value = 33402.855
$ws.Range("L83").Value = 100208.565
$ws.Range("N83").Value = -110192.565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 22897.6
$ws.Range("J35").Value = 24872
$ws.Range("L35").Value = 24872
$ws.Range("N35").Value = -25492
$ws.Range("H62").Value = 92496.5
$ws.Range("J62").Value = 92496.5
$ws.Range("L62").Value = 92496.5
$ws.Range("N62").Value = -93868.5
$ws.Range("H65").Value = 92496.5
$ws.Range("J65").Value = 92496.5
$ws.Range("L65").Value = 277489.5
$ws.Range("N65").Value = -284353.5
$ws.Range("H82").Value = 16844.125
$ws.Range("I82").Value = 2993.8572
$ws.Range("J82").Value = 27616.555
$ws.Range("K82").Value = 2993.8572
$ws.Range("L82").Value = 27616.555
$ws.Range("M82").Value = -2610.8572
$ws.Range("N82").Value = -28382.555
$ws.Range("H85").Value = 16844.125
$ws.Range("I85").Value = 2993.8572
$ws.Range("J85").Value = 27616.555
$ws.Range("K85").Value = 2993.8572
$ws.Range("L85").Value = 27616.555
$ws.Range("M85").Value = -1667.8572
$ws.Range("N85").Value = -30268.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16395403
$ws.Range("I31").Value = 43479676
$ws.Range("J31").Value = 2289.5
$ws.Range("K31").Value = 43479676
$ws.Range("L31").Value = 2289.5
$ws.Range("M31").Value = -43479381
$ws.Range("N31").Value = -2879.5
$ws.Range("H34").Value = 16395403
$ws.Range("I34").Value = 43479676
$ws.Range("J34").Value = 2289.5
$ws.Range("K34").Value = 43479676
$ws.Range("L34").Value = 2289.5
$ws.Range("M34").Value = -43479474
$ws.Range("N34").Value = -2693.5
$ws.Range("H41").Value = 19241.166
$ws.Range("J41").Value = 22077.6
$ws.Range("L41").Value = 22077.6
$ws.Range("N41").Value = -22933.6
$ws.Range("H45").Value = 10511.167
$ws.Range("I45").Value = 3067
$ws.Range("J45").Value = 12000
$ws.Range("K45").Value = 3067
$ws.Range("L45").Value = 12000
$ws.Range("M45").Value = -2474
$ws.Range("N45").Value = -13186
$ws.Range("H50").Value = 14396
$ws.Range("J50").Value = 14396
$ws.Range("L50").Value = 14396
$ws.Range("N50").Value = -15646
$ws.Range("H51").Value = 25726.1
$ws.Range("J51").Value = 28142.334
$ws.Range("L51").Value = 28142.334
$ws.Range("N51").Value = -29614.334
$ws.Range("H60").Value = 10854.333
$ws.Range("J60").Value = 14500.6
$ws.Range("L60").Value = 14500.6
$ws.Range("N60").Value = -15522.6
$ws.Range("H61").Value = 25726.1
$ws.Range("J61").Value = 28142.334
$ws.Range("L61").Value = 28142.334
$ws.Range("N61").Value = -28838.334
$ws.Range("H63").Value = 33000
$ws.Range("J63").Value = 33000
$ws.Range("L63").Value = 33000
$ws.Range("N63").Value = -34372
$ws.Range("H66").Value = 33000
$ws.Range("J66").Value = 33000
$ws.Range("L66").Value = 99000
$ws.Range("N66").Value = -105864
$ws.Range("H68").Value = 25984.285
$ws.Range("J68").Value = 25984.285
$ws.Range("L68").Value = 25984.285
$ws.Range("N68").Value = -27482.285
$ws.Range("H71").Value = 25984.285
$ws.Range("J71").Value = 25984.285
$ws.Range("L71").Value = 77952.855
$ws.Range("N71").Value = -85440.855
$ws.Range("H109").Value = 27200
$ws.Range("J109").Value = 27200
$ws.Range("L109").Value = 27200
$ws.Range("N109").Value = -29280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2607.5386
$ws.Range("I51").Value = 426
$ws.Range("J51").Value = 3577.111
$ws.Range("K51").Value = 1278
$ws.Range("L51").Value = 10731.333
$ws.Range("M51").Value = -818
$ws.Range("N51").Value = -11651.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 26404.5
$ws.Range("J57").Value = 26404.5
$ws.Range("L57").Value = 26404.5
$ws.Range("N57").Value = -28044.5
$ws.Range("H64").Value = 34444
$ws.Range("J64").Value = 34444
$ws.Range("L64").Value = 34444
$ws.Range("N64").Value = -34940
$ws.Range("H67").Value = 34444
$ws.Range("J67").Value = 34444
$ws.Range("L67").Value = 34444
$ws.Range("N67").Value = -36160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1864
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 1478.5
$ws.Range("I46").Value = 1322
$ws.Range("K46").Value = 1322
$ws.Range("M46").Value = -1134
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26622
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -83112
$ws.Range("H122").Value = 4033.9285
$ws.Range("I122").Value = 3816.6667
$ws.Range("J122").Value = 4196.875
$ws.Range("K122").Value = 11450.0001
$ws.Range("L122").Value = 12590.625
$ws.Range("M122").Value = -9000.000100000001
$ws.Range("N122").Value = -17490.625
$ws.Range("H132").Value = 2427.1428
$ws.Range("I132").Value = 1789.2727
$ws.Range("J132").Value = 4766
$ws.Range("K132").Value = 5367.8181
$ws.Range("L132").Value = 14298
$ws.Range("M132").Value = -2837.8181
$ws.Range("N132").Value = -19358

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12997.333
$ws.Range("I62").Value = 9002
$ws.Range("J62").Value = 14995
$ws.Range("K62").Value = 9002
$ws.Range("L62").Value = 14995
$ws.Range("M62").Value = -8378
$ws.Range("N62").Value = -16243
$ws.Range("H65").Value = 12997.333
$ws.Range("I65").Value = 9002
$ws.Range("J65").Value = 14995
$ws.Range("K65").Value = 45010
$ws.Range("L65").Value = 74975
$ws.Range("M65").Value = -41890
$ws.Range("N65").Value = -81215

